$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 2644.625
$ws.Range("I2").Value = 3272
$ws.Range("J2").Value = 1599
$ws.Range("K2").Value = 3272
$ws.Range("L2").Value = 1599
$ws.Range("M2").Value = -3159
$ws.Range("N2").Value = -1825
$ws.Range("H5").Value = 102.888885
$ws.Range("I5").Value = 102.888885
$ws.Range("K5").Value = 102.888885
$ws.Range("M5").Value = 12.111115
$ws.Range("H18").Value = 6559.8
$ws.Range("J18").Value = 12749.5
$ws.Range("L18").Value = 12749.5
$ws.Range("N18").Value = -13317.5
$ws.Range("H28").Value = 976.26666
$ws.Range("I28").Value = 1022.5455
$ws.Range("K28").Value = 1022.5455
$ws.Range("M28").Value = -537.5454999999999
$ws.Range("H33").Value = 125.818184
$ws.Range("I33").Value = 114.25
$ws.Range("K33").Value = 114.25
$ws.Range("M33").Value = 114.75
$ws.Range("H40").Value = 2380
$ws.Range("I40").Value = 2380
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 2380
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -2205
$ws.Range("N40").ClearContents()
$ws.Range("H42").Value = 587.3
$ws.Range("I42").Value = 624.1429000000001
$ws.Range("J42").Value = 501.33334
$ws.Range("K42").Value = 1872.4287
$ws.Range("L42").Value = 1504.00002
$ws.Range("M42").Value = -1642.4287
$ws.Range("N42").Value = -1964.00002
$ws.Range("H53").Value = 750
$ws.Range("J53").Value = 833.3333
$ws.Range("L53").Value = 833.3333
$ws.Range("N53").Value = -2107.3333
$ws.Range("H64").Value = 72148.664
$ws.Range("J64").Value = 10000
$ws.Range("L64").Value = 10000
$ws.Range("N64").Value = -10496
$ws.Range("H67").Value = 72148.664
$ws.Range("J67").Value = 10000
$ws.Range("L67").Value = 10000
$ws.Range("N67").Value = -11716
$ws.Range("H70").Value = 13899031
$ws.Range("I70").Value = 22223848
$ws.Range("J70").Value = 24333.334
$ws.Range("K70").Value = 66671544
$ws.Range("L70").Value = 73000.00199999999
$ws.Range("M70").Value = -66671274
$ws.Range("N70").Value = -73540.00199999999
$ws.Range("H73").Value = 13899031
$ws.Range("I73").Value = 22223848
$ws.Range("J73").Value = 24333.334
$ws.Range("K73").Value = 66671544
$ws.Range("L73").Value = 73000.00199999999
$ws.Range("M73").Value = -66670608
$ws.Range("N73").Value = -74872.00199999999
$ws.Range("H76").Value = 2939.625
$ws.Range("I76").Value = 2645.2856
$ws.Range("K76").Value = 2645.2856
$ws.Range("M76").Value = -2330.2856
$ws.Range("H79").Value = 2939.625
$ws.Range("I79").Value = 2645.2856
$ws.Range("K79").Value = 2645.2856
$ws.Range("M79").Value = -1553.2856
$ws.Range("H80").Value = 1557699.6
$ws.Range("J80").Value = 7239.8335
$ws.Range("L80").Value = 21719.5005
$ws.Range("N80").Value = -23715.5005
$ws.Range("H83").Value = 1557699.6
$ws.Range("J83").Value = 7239.8335
$ws.Range("L83").Value = 65158.5015
$ws.Range("N83").Value = -75142.5015
$ws.Range("H88").Value = 6862.1816
$ws.Range("I88").Value = 3948.8
$ws.Range("J88").Value = 9290
$ws.Range("K88").Value = 3948.8
$ws.Range("L88").Value = 9290
$ws.Range("M88").Value = -3542.8
$ws.Range("N88").Value = -10102
$ws.Range("H91").Value = 6862.1816
$ws.Range("I91").Value = 3948.8
$ws.Range("J91").Value = 9290
$ws.Range("K91").Value = 3948.8
$ws.Range("L91").Value = 9290
$ws.Range("M91").Value = -2544.8
$ws.Range("N91").Value = -12098
$ws.Range("H106").Value = 2242.0527
$ws.Range("I106").Value = 1959.4706
$ws.Range("K106").Value = 1959.4706
$ws.Range("M106").Value = -1328.4706
$ws.Range("H113").Value = 6058.4
$ws.Range("J113").Value = 6159.857
$ws.Range("L113").Value = 6159.857
$ws.Range("N113").Value = -12667.857
$ws.Range("H116").Value = 10708.777
$ws.Range("I116").Value = 12000.963
$ws.Range("J116").Value = 6832.222
$ws.Range("K116").Value = 12000.963
$ws.Range("L116").Value = 6832.222
$ws.Range("M116").Value = -8558.963
$ws.Range("N116").Value = -13716.222
$ws.Range("H125").Value = 6225903
$ws.Range("I125").Value = 20452210
$ws.Range("J125").Value = 1893.1875
$ws.Range("K125").Value = 184069890
$ws.Range("L125").Value = 17038.6875
$ws.Range("M125").Value = -184067430
$ws.Range("N125").Value = -21958.6875
$ws.Range("H127").Value = 1420.3334
$ws.Range("I127").Value = 752.5714
$ws.Range("K127").Value = 2257.7142
$ws.Range("M127").Value = 2702.2858
$ws.Range("H132").Value = 3760642
$ws.Range("I132").Value = 3862271.2
$ws.Range("J132").Value = 365
$ws.Range("K132").Value = 11586813.6
$ws.Range("L132").Value = 1095
$ws.Range("M132").Value = -11584283.6
$ws.Range("N132").Value = -6155
$ws.Range("H135").Value = 2286.1304
$ws.Range("I135").Value = 1904.7894
$ws.Range("J135").Value = 4097.5
$ws.Range("K135").Value = 17143.1046
$ws.Range("L135").Value = 36877.5
$ws.Range("M135").Value = -14608.1046
$ws.Range("N135").Value = -41947.5
$ws.Range("H137").Value = 8224.319
$ws.Range("I137").Value = 1210.1052
$ws.Range("J137").Value = 12983.964
$ws.Range("K137").Value = 3630.3156
$ws.Range("L137").Value = 38951.892
$ws.Range("M137").Value = -1080.3156
$ws.Range("N137").Value = -44051.892
$ws.Range("H138").Value = 1220569.9
$ws.Range("I138").Value = 2603.0667
$ws.Range("J138").Value = 1593416.9
$ws.Range("K138").Value = 7809.2001
$ws.Range("L138").Value = 4780250.699999999
$ws.Range("M138").Value = -2669.2001
$ws.Range("N138").Value = -4790530.699999999
$ws.Range("H141").Value = 2355.2856
$ws.Range("I141").Value = 2355.2856
$ws.Range("K141").Value = 7065.8568
$ws.Range("M141").Value = -1885.8568

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H29").Value = 2466.6667
$ws.Range("I29").Value = 2466.6667
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 2466.6667
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -2158.6667
$ws.Range("N29").ClearContents()
$ws.Range("H32").Value = 18271.904
$ws.Range("I32").Value = 18872.6
$ws.Range("J32").Value = 6258
$ws.Range("K32").Value = 18872.6
$ws.Range("L32").Value = 6258
$ws.Range("M32").Value = -18585.6
$ws.Range("N32").Value = -6832
$ws.Range("H45").Value = 5192.222
$ws.Range("I45").Value = 4524.1
$ws.Range("J45").Value = 6027.375
$ws.Range("K45").Value = 4524.1
$ws.Range("L45").Value = 6027.375
$ws.Range("M45").Value = -4147.1
$ws.Range("N45").Value = -6781.375
$ws.Range("H61").Value = 6512.75
$ws.Range("I61").Value = 3832.4644
$ws.Range("J61").Value = 15893.75
$ws.Range("K61").Value = 3832.4644
$ws.Range("L61").Value = 15893.75
$ws.Range("M61").Value = -3620.4644
$ws.Range("N61").Value = -16317.75
$ws.Range("H74").Value = 3192.5862
$ws.Range("I74").Value = 1192.238
$ws.Range("J74").Value = 8443.5
$ws.Range("K74").Value = 1192.238
$ws.Range("L74").Value = 8443.5
$ws.Range("M74").Value = -318.2380000000001
$ws.Range("N74").Value = -10191.5
$ws.Range("H77").Value = 3192.5862
$ws.Range("I77").Value = 1192.238
$ws.Range("J77").Value = 8443.5
$ws.Range("K77").Value = 5961.190000000001
$ws.Range("L77").Value = 42217.5
$ws.Range("M77").Value = -1593.190000000001
$ws.Range("N77").Value = -50953.5
$ws.Range("H88").Value = 17638
$ws.Range("J88").Value = 17638
$ws.Range("L88").Value = 17638
$ws.Range("N88").Value = -18450
$ws.Range("H91").Value = 17638
$ws.Range("J91").Value = 17638
$ws.Range("L91").Value = 17638
$ws.Range("N91").Value = -20446
$ws.Range("H122").Value = 1653.6875
$ws.Range("I122").Value = 1323.1428
$ws.Range("J122").Value = 3967.5
$ws.Range("K122").Value = 3969.4284
$ws.Range("L122").Value = 11902.5
$ws.Range("M122").Value = -1519.4284
$ws.Range("N122").Value = -16802.5
$ws.Range("H132").Value = 1734.5238
$ws.Range("I132").Value = 1418.4117
$ws.Range("J132").Value = 3078
$ws.Range("K132").Value = 4255.2351
$ws.Range("L132").Value = 9234
$ws.Range("M132").Value = -1725.2351
$ws.Range("N132").Value = -14294
$ws.Range("H136").Value = 6512.75
$ws.Range("I136").Value = 3832.4644
$ws.Range("J136").Value = 15893.75
$ws.Range("K136").Value = 11497.3932
$ws.Range("L136").Value = 47681.25
$ws.Range("M136").Value = -8947.393199999999
$ws.Range("N136").Value = -52781.25

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2858.3333
$ws.Range("I86").Value = 3366.6667
$ws.Range("J86").Value = 2350
$ws.Range("K86").Value = 3366.6667
$ws.Range("L86").Value = 2350
$ws.Range("M86").Value = -2243.6667
$ws.Range("N86").Value = -4596
$ws.Range("H89").Value = 2858.3333
$ws.Range("I89").Value = 3366.6667
$ws.Range("J89").Value = 2350
$ws.Range("K89").Value = 16833.3335
$ws.Range("L89").Value = 11750
$ws.Range("M89").Value = -11217.3335
$ws.Range("N89").Value = -22982
$ws.Range("H94").Value = 5000779
$ws.Range("J94").Value = 14286184
$ws.Range("L94").Value = 14286184
$ws.Range("N94").Value = -14287086
$ws.Range("H99").Value = 1599.6666
$ws.Range("I99").Value = 1219.4286
$ws.Range("J99").Value = 2930.5
$ws.Range("K99").Value = 1219.4286
$ws.Range("L99").Value = 2930.5
$ws.Range("M99").Value = 278.5714
$ws.Range("N99").Value = -5926.5
$ws.Range("H105").Value = 2262.7097
$ws.Range("I105").Value = 2047.1578
$ws.Range("J105").Value = 2604
$ws.Range("K105").Value = 2047.1578
$ws.Range("L105").Value = 2604
$ws.Range("M105").Value = -300.1578
$ws.Range("N105").Value = -6098
$ws.Range("H134").Value = 9562.297
$ws.Range("I134").Value = 10232.971
$ws.Range("J134").Value = 1961.3334
$ws.Range("K134").Value = 30698.913
$ws.Range("L134").Value = 5884.0002
$ws.Range("M134").Value = -28163.913
$ws.Range("N134").Value = -10954.0002
$ws.Range("H137").Value = 54499.75
$ws.Range("I137").Value = 70000
$ws.Range("K137").Value = 70000
$ws.Range("M137").Value = -64900

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 10614832
$ws.Range("I6").Value = 15921248
$ws.Range("K6").Value = 15921248
$ws.Range("M6").Value = -15921135
$ws.Range("H16").Value = 3704.5715
$ws.Range("I16").Value = 3888
$ws.Range("J16").Value = 3246
$ws.Range("K16").Value = 3888
$ws.Range("L16").Value = 3246
$ws.Range("M16").Value = -3601
$ws.Range("N16").Value = -3820
$ws.Range("H31").Value = 1963082.6
$ws.Range("I31").Value = 5002446
$ws.Range("J31").Value = 2203.0322
$ws.Range("K31").Value = 5002446
$ws.Range("L31").Value = 2203.0322
$ws.Range("M31").Value = -5002151
$ws.Range("N31").Value = -2793.0322
$ws.Range("H34").Value = 1963082.6
$ws.Range("I34").Value = 5002446
$ws.Range("J34").Value = 2203.0322
$ws.Range("K34").Value = 5002446
$ws.Range("L34").Value = 2203.0322
$ws.Range("M34").Value = -5002244
$ws.Range("N34").Value = -2607.0322
$ws.Range("H41").Value = 40000
$ws.Range("J41").Value = 40000
$ws.Range("L41").Value = 40000
$ws.Range("N41").Value = -40856
$ws.Range("H50").Value = 36995
$ws.Range("J50").Value = 36995
$ws.Range("L50").Value = 36995
$ws.Range("N50").Value = -38245
$ws.Range("H58").Value = 1397.7812
$ws.Range("I58").Value = 919.3
$ws.Range("K58").Value = 919.3
$ws.Range("M58").Value = -716.3
$ws.Range("H76").Value = 5000
$ws.Range("I76").Value = 5000
$ws.Range("K76").Value = 5000
$ws.Range("M76").Value = -4685
$ws.Range("H79").Value = 5000
$ws.Range("I79").Value = 5000
$ws.Range("K79").Value = 5000
$ws.Range("M79").Value = -3908
$ws.Range("H99").Value = 9219.556
$ws.Range("I99").Value = 10496.75
$ws.Range("J99").Value = 8197.799999999999
$ws.Range("K99").Value = 10496.75
$ws.Range("L99").Value = 8197.799999999999
$ws.Range("M99").Value = -8998.75
$ws.Range("N99").Value = -11193.8
$ws.Range("H105").Value = 1788
$ws.Range("I105").Value = 717.3333
$ws.Range("K105").Value = 717.3333
$ws.Range("M105").Value = 1029.6667
$ws.Range("H113").Value = 3704.5715
$ws.Range("I113").Value = 3888
$ws.Range("J113").Value = 3246
$ws.Range("K113").Value = 3888
$ws.Range("L113").Value = 3246
$ws.Range("M113").Value = -1718
$ws.Range("N113").Value = -7586
$ws.Range("H126").Value = 9219.556
$ws.Range("I126").Value = 10496.75
$ws.Range("J126").Value = 8197.799999999999
$ws.Range("K126").Value = 31490.25
$ws.Range("L126").Value = 24593.4
$ws.Range("M126").Value = -29020.25
$ws.Range("N126").Value = -29533.4
$ws.Range("H132").Value = 3768.889
$ws.Range("I132").Value = 2928
$ws.Range("J132").Value = 4441.6
$ws.Range("K132").Value = 8784
$ws.Range("L132").Value = 13324.8
$ws.Range("M132").Value = -6254
$ws.Range("N132").Value = -18384.8
$ws.Range("H134").Value = 1819.2916
$ws.Range("I134").Value = 1269.0975
$ws.Range("J134").Value = 5041.857
$ws.Range("K134").Value = 3807.2925
$ws.Range("L134").Value = 15125.571
$ws.Range("M134").Value = -1272.2925
$ws.Range("N134").Value = -20195.571
$ws.Range("H136").Value = 1397.7812
$ws.Range("I136").Value = 919.3
$ws.Range("K136").Value = 2757.9
$ws.Range("M136").Value = -207.8999999999996

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 752.9231
$ws.Range("I5").Value = 683
$ws.Range("J5").Value = 812.8570999999999
$ws.Range("K5").Value = 2049
$ws.Range("L5").Value = 2438.5713
$ws.Range("M5").Value = -1937
$ws.Range("N5").Value = -2662.5713
$ws.Range("H61").Value = 70.625
$ws.Range("I61").Value = 90
$ws.Range("J61").Value = 38.333332
$ws.Range("K61").Value = 270
$ws.Range("L61").Value = 114.999996
$ws.Range("M61").Value = -55
$ws.Range("N61").Value = -544.999996
$ws.Range("H105").Value = 7595.5
$ws.Range("J105").Value = 7962.5713
$ws.Range("L105").Value = 23887.7139
$ws.Range("N105").Value = -29129.7139
$ws.Range("H107").Value = 1512.3972
$ws.Range("I107").Value = 708.25806
$ws.Range("J107").Value = 2105.9285
$ws.Range("K107").Value = 2124.77418
$ws.Range("L107").Value = 6317.7855
$ws.Range("M107").Value = -204.7741799999999
$ws.Range("N107").Value = -10157.7855
$ws.Range("H122").Value = 667.2273
$ws.Range("J122").Value = 705.44446
$ws.Range("L122").Value = 6349.00014
$ws.Range("N122").Value = -11249.00014
$ws.Range("H132").Value = 1625.8695
$ws.Range("I132").Value = 1558.875
$ws.Range("J132").Value = 1661.6
$ws.Range("K132").Value = 14029.875
$ws.Range("L132").Value = 14954.4
$ws.Range("M132").Value = -11499.875
$ws.Range("N132").Value = -20014.4
$ws.Range("H135").Value = 752.9231
$ws.Range("I135").Value = 683
$ws.Range("J135").Value = 812.8570999999999
$ws.Range("K135").Value = 6147
$ws.Range("L135").Value = 7315.7139
$ws.Range("M135").Value = -3612
$ws.Range("N135").Value = -12385.7139

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3603.6667
$ws.Range("I80").Value = 2781
$ws.Range("J80").Value = 5249
$ws.Range("K80").Value = 2781
$ws.Range("L80").Value = 5249
$ws.Range("M80").Value = -1783
$ws.Range("N80").Value = -7245
$ws.Range("H83").Value = 3603.6667
$ws.Range("I83").Value = 2781
$ws.Range("J83").Value = 5249
$ws.Range("K83").Value = 13905
$ws.Range("L83").Value = 26245
$ws.Range("M83").Value = -8913
$ws.Range("N83").Value = -36229
$ws.Range("H100").Value = 75000
$ws.Range("J100").Value = 75000
$ws.Range("L100").Value = 75000
$ws.Range("N100").Value = -77164
$ws.Range("H113").Value = 1668.1818
$ws.Range("I113").Value = 992.4286
$ws.Range("K113").Value = 992.4286
$ws.Range("M113").Value = 1177.5714
$ws.Range("H132").Value = 1674.6471
$ws.Range("I132").Value = 1527.5778
$ws.Range("J132").Value = 2777.6667
$ws.Range("K132").Value = 4582.7334
$ws.Range("L132").Value = 8333.000100000001
$ws.Range("M132").Value = -2052.7334
$ws.Range("N132").Value = -13393.0001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2197.5
$ws.Range("I22").Value = 2190
$ws.Range("K22").Value = 2190
$ws.Range("M22").Value = -1895
$ws.Range("H27").Value = 2197.5
$ws.Range("I27").Value = 2190
$ws.Range("K27").Value = 2190
$ws.Range("M27").Value = -2083
$ws.Range("H31").Value = 6395.6665
$ws.Range("I31").Value = 1996.3334
$ws.Range("J31").Value = 10795
$ws.Range("K31").Value = 1996.3334
$ws.Range("L31").Value = 10795
$ws.Range("M31").Value = -1748.3334
$ws.Range("N31").Value = -11291
$ws.Range("H34").Value = 9999.5
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()
$ws.Range("H46").Value = 4075.6365
$ws.Range("I46").Value = 999.5
$ws.Range("K46").Value = 999.5
$ws.Range("M46").Value = -811.5
$ws.Range("H68").Value = 1825.2222
$ws.Range("I68").Value = 1886
$ws.Range("K68").Value = 1886
$ws.Range("M68").Value = -1137
$ws.Range("H71").Value = 1825.2222
$ws.Range("I71").Value = 1886
$ws.Range("K71").Value = 9430
$ws.Range("M71").Value = -5686
$ws.Range("H82").Value = 3446
$ws.Range("I82").Value = 2910.8462
$ws.Range("J82").Value = 3981.1538
$ws.Range("K82").Value = 2910.8462
$ws.Range("L82").Value = 3981.1538
$ws.Range("M82").Value = -2549.8462
$ws.Range("N82").Value = -4703.1538
$ws.Range("H85").Value = 3446
$ws.Range("I85").Value = 2910.8462
$ws.Range("J85").Value = 3981.1538
$ws.Range("K85").Value = 2910.8462
$ws.Range("L85").Value = 3981.1538
$ws.Range("M85").Value = -1662.8462
$ws.Range("N85").Value = -6477.1538
$ws.Range("H122").Value = 4116.6665
$ws.Range("I122").Value = 3400
$ws.Range("J122").Value = 4833.3335
$ws.Range("K122").Value = 10200
$ws.Range("L122").Value = 14500.0005
$ws.Range("M122").Value = -7750
$ws.Range("N122").Value = -19400.0005
$ws.Range("H132").Value = 3770.9807
$ws.Range("I132").Value = 2983.3125
$ws.Range("J132").Value = 5031.25
$ws.Range("K132").Value = 8949.9375
$ws.Range("L132").Value = 15093.75
$ws.Range("M132").Value = -6419.9375
$ws.Range("N132").Value = -20153.75
$ws.Range("H136").Value = 3122.3635
$ws.Range("I136").Value = 1634.5714
$ws.Range("K136").Value = 4903.7142
$ws.Range("M136").Value = -2353.7142

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 29969
$ws.Range("J20").Value = 29969
$ws.Range("L20").Value = 29969
$ws.Range("N20").Value = -30449
$ws.Range("H81").Value = 5096.9688
$ws.Range("I81").Value = 4957.8076
$ws.Range("K81").Value = 9915.6152
$ws.Range("M81").Value = -8854.6152
$ws.Range("H84").Value = 5096.9688
$ws.Range("I84").Value = 4957.8076
$ws.Range("K84").Value = 49578.076
$ws.Range("M84").Value = -44274.076
$ws.Range("H107").Value = 807.2778
$ws.Range("I107").Value = 849.0714
$ws.Range("K107").Value = 2547.2142
$ws.Range("M107").Value = -627.2142000000003
$ws.Range("H113").Value = 650.27905
$ws.Range("I113").Value = 694.02856
$ws.Range("J113").Value = 458.875
$ws.Range("K113").Value = 2082.08568
$ws.Range("L113").Value = 1376.625
$ws.Range("M113").Value = 87.91431999999986
$ws.Range("N113").Value = -5716.625
$ws.Range("H122").Value = 5850.6665
$ws.Range("I122").Value = 6662.2354
$ws.Range("J122").Value = 3879.7144
$ws.Range("K122").Value = 19986.7062
$ws.Range("L122").Value = 11639.1432
$ws.Range("M122").Value = -17536.7062
$ws.Range("N122").Value = -16539.1432
$ws.Range("H126").Value = 3944.8462
$ws.Range("I126").Value = 2950.75
$ws.Range("J126").Value = 4386.6665
$ws.Range("K126").Value = 8852.25
$ws.Range("L126").Value = 13159.9995
$ws.Range("M126").Value = -6382.25
$ws.Range("N126").Value = -18099.9995
$ws.Range("H132").Value = 14195.258
$ws.Range("I132").Value = 15025.705
$ws.Range("J132").Value = 4063.8
$ws.Range("K132").Value = 45077.115
$ws.Range("L132").Value = 12191.4
$ws.Range("M132").Value = -42547.115
$ws.Range("N132").Value = -17251.4
$ws.Range("H136").Value = 6868
$ws.Range("I136").Value = 7915.8823
$ws.Range("K136").Value = 23747.6469
$ws.Range("M136").Value = -21197.6469
